$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-4 with new values
$ws.Range("A2").Value = 22
$ws.Range("B2").Value = 2

$ws.Range("A3").Value = 21
$ws.Range("B3").Value = 2

$ws.Range("A4").Value = 11
$ws.Range("B4").Value = 2

# Add new row 5
$ws.Range("A5").Value = 12
$ws.Range("B5").Value = 1

# Copy style from A4 (which carries style index 1) to the new A5 cell
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)  # xlPasteFormats
